# Regenerate save_data to use K instead of Strike#, write new s_vals values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
